$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timestamp for all data rows (2-38)
$newDate = "2023-05-25 13:07:26"
for ($r = 2; $r -le 38; $r++) {
    $ws.Range("B$r").Value = $newDate
}

# Row 2
$ws.Range("E2").Value = 60564

# Row 3
$ws.Range("E3").Value = 4256

# Row 5
$ws.Range("E5").Value = 390886
$ws.Range("F5").Value = 111

# Row 6
$ws.Range("E6").Value = 4733942
$ws.Range("F6").Value = 10

# Row 7
$ws.Range("E7").Value = 42553
$ws.Range("F7").Value = 90

# Row 8
$ws.Range("E8").Value = 319482

# Row 9
$ws.Range("E9").Value = 2412896

# Row 10
$ws.Range("E10").Value = 42734
$ws.Range("F10").Value = 83

# Row 11
$ws.Range("E11").Value = 23993
$ws.Range("F11").Value = 76

# Row 12
$ws.Range("E12").Value = 1973401
$ws.Range("F12").Value = 9

# Row 13
$ws.Range("E13").Value = 1047459

# Row 14
$ws.Range("E14").Value = 256189
$ws.Range("F14").Value = 105

# Row 15
$ws.Range("E15").Value = 114733
$ws.Range("F15").Value = 174

# Row 16
$ws.Range("E16").Value = 72015
$ws.Range("F16").Value = 118

# Row 17
$ws.Range("E17").Value = 170

# Row 18
$ws.Range("E18").Value = 405414
$ws.Range("F18").Value = 54

# Row 19
$ws.Range("E19").Value = 3991605
$ws.Range("F19").Value = 12

# Row 20
$ws.Range("E20").Value = 4103
$ws.Range("F20").ClearContents()

# Row 21
$ws.Range("E21").Value = 6525

# Row 22
$ws.Range("E22").Value = 3091
$ws.Range("F22").Value = 58

# Row 23
$ws.Range("E23").Value = 45091
$ws.Range("F23").Value = 178

# Row 24
$ws.Range("E24").Value = 860618
$ws.Range("F24").Value = 7

# Row 25
$ws.Range("E25").Value = 206

# Row 26
$ws.Range("E26").Value = 1563

# Row 27
$ws.Range("D27").Value = 3.7
$ws.Range("E27").Value = 14

# Row 29
$ws.Range("E29").Value = 306

# Row 31
$ws.Range("D31").Value = 1.3
$ws.Range("E31").Value = 24

# Row 32
$ws.Range("E32").Value = 1362

# Row 33
$ws.Range("E33").Value = 27433

# Row 34
$ws.Range("C34").Value = "App Store"
$ws.Range("D34").ClearContents()
$ws.Range("E34").ClearContents()

# Row 35
$ws.Range("E35").Value = 13130

# Row 36
$ws.Range("E36").Value = 35

# Row 37
$ws.Range("E37").Value = 3215

Write-Host "Applied iOS and Android Fix updates"